$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.958190381526947
$ws.Range("B1").Value = 1.771008849143982
$ws.Range("C1").Value = 6.801307201385498
$ws.Range("D1").Value = 1.61429750919342
$ws.Range("E1").Value = 0.9321849346160889
